# Apply the grade-schedule text corrections described by the commit
# "add: creating objective var to direct search".
#
# The workbook has two sheets; the data that changed lives on the first
# sheet ("Cola aqui os valores"), which is not necessarily the
# ActiveSheet, so we look it up explicitly by name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cola aqui os valores")

# Row 3
$ws.Range("B3").Value = "Inglês"
$ws.Range("C3").Value = "Inglês"
$ws.Range("D3").Value = "História"
$ws.Range("E3").Value = "Geografia"
$ws.Range("I3").Value = "Matemática"
$ws.Range("K3").Value = "Língua Portuguesa"
$ws.Range("L3").Value = "Inglês"
$ws.Range("M3").Value = "Matemática"

# Row 4
$ws.Range("B4").Value = "Mind Makers"
$ws.Range("C4").Value = "Geografia"
$ws.Range("E4").Value = "Língua Portuguesa"
$ws.Range("F4").Value = "Inglês"
$ws.Range("I4").Value = "Inglês"
$ws.Range("L4").Value = "História"
$ws.Range("M4").Value = "Matemática"

# Row 5
$ws.Range("C5").Value = "Língua Portuguesa"
$ws.Range("D5").Value = "Inglês"
$ws.Range("F5").Value = "Língua Portuguesa"
$ws.Range("K5").Value = "Geografia"
$ws.Range("L5").Value = "Língua Portuguesa"
$ws.Range("M5").Value = "Inglês"

# Row 7
$ws.Range("B7").Value = "Matemática"
$ws.Range("F7").Value = "Língua Portuguesa"
$ws.Range("I7").Value = "Língua Portuguesa"
$ws.Range("J7").Value = "Ciências"
$ws.Range("M7").Value = "Música"

# Row 8
$ws.Range("B8").Value = "Música"
$ws.Range("F8").Value = "Educação Física"
$ws.Range("I8").Value = "Matemática"
$ws.Range("J8").Value = "Ciências"
$ws.Range("L8").Value = "Matemática"

# Row 13
$ws.Range("B13").Value = "Artes"
$ws.Range("C13").Value = "Língua Portuguesa"
$ws.Range("D13").Value = "Mind Makers"
$ws.Range("E13").Value = "Língua Portuguesa"
$ws.Range("F13").Value = "Artes"
$ws.Range("I13").Value = "História"
$ws.Range("J13").Value = "Ed. Financeira"
$ws.Range("K13").Value = "Inglês"
$ws.Range("L13").Value = "Matemática"
$ws.Range("M13").Value = "Inglês"

# Row 14
$ws.Range("B14").Value = "Geografia"
$ws.Range("C14").Value = "Língua Portuguesa"
$ws.Range("D14").Value = "História"
$ws.Range("E14").Value = "Inglês"
$ws.Range("F14").Value = "Música"
$ws.Range("I14").Value = "Língua Portuguesa"
$ws.Range("J14").Value = "Inglês"
$ws.Range("K14").Value = "Matemática"
$ws.Range("L14").Value = "Matemática"
$ws.Range("M14").Value = "Língua Portuguesa"

# Row 15
$ws.Range("B15").Value = "Geografia"
$ws.Range("C15").Value = "Inglês"
$ws.Range("D15").Value = "Ed. Financeira"
$ws.Range("E15").Value = "História"
$ws.Range("F15").Value = "Educação Física"
$ws.Range("J15").Value = "Geografia"
$ws.Range("K15").Value = "Língua Portuguesa"
$ws.Range("L15").Value = "Inglês"

# Row 17
$ws.Range("B17").Value = "Ciências"
$ws.Range("C17").Value = "Matemática"
$ws.Range("F17").Value = "Matemática"
$ws.Range("K17").Value = "Mind Makers"
$ws.Range("M17").Value = "Educação Física"

# Row 18
$ws.Range("B18").Value = "Ciências"
$ws.Range("C18").Value = "Matemática"
$ws.Range("E18").Value = "Língua Portuguesa"
$ws.Range("F18").Value = "Matemática"
$ws.Range("I18").Value = "Língua Portuguesa"
$ws.Range("M18").Value = "Música"
